# Scheduled data refresh: update market-price / profit columns (H-N)
# across the Ultima_Profits crafting-profit sheets (ALC, ARM, CRP, CUL, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ALC row 94
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H94").Value = 6000
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 6000
$ws.Range("K94").Value = 0
$ws.Range("L94").ClearContents()
$ws.Range("M94").Value = 6000
$ws.Range("N94").Value = -6902

# ALC row 96
$ws.Range("H96").Value = 1680.7
$ws.Range("I96").Value = 1096.1666
$ws.Range("J96").Value = 2557.5
$ws.Range("K96").Value = 3288.4998
$ws.Range("L96").Value = 7672.5
$ws.Range("M96").Value = -1915.4998
$ws.Range("N96").Value = -10418.5

# ALC row 112
$ws.Range("H112").Value = 1951.5358
$ws.Range("I112").Value = 1199.5
$ws.Range("J112").Value = 2009.3846
$ws.Range("K112").Value = 3598.5
$ws.Range("L112").Value = 6028.1538
$ws.Range("M112").Value = -2490.5
$ws.Range("N112").Value = -8244.1538

# ALC row 113
$ws.Range("H113").Value = 46648.125
$ws.Range("I113").Value = 59089.44
$ws.Range("J113").Value = 2214.8572
$ws.Range("K113").Value = 59089.44
$ws.Range("L113").Value = 2214.8572
$ws.Range("M113").Value = -55835.44
$ws.Range("N113").Value = -8722.8572

# ALC row 125
$ws.Range("H125").Value = 1388.091
$ws.Range("I125").Value = 2335.8
$ws.Range("K125").Value = 21022.2
$ws.Range("M125").Value = -18562.2

# ALC row 132
$ws.Range("H132").Value = 4646.4424
$ws.Range("I132").Value = 4327.6514
$ws.Range("J132").Value = 6169.5557
$ws.Range("K132").Value = 12982.9542
$ws.Range("L132").Value = 18508.6671
$ws.Range("M132").Value = -10452.9542
$ws.Range("N132").Value = -23568.6671

# ALC row 138
$ws.Range("H138").Value = 1812.78
$ws.Range("I138").Value = 805
$ws.Range("J138").Value = 2513.1018
$ws.Range("K138").Value = 2415
$ws.Range("L138").Value = 7539.305399999999
$ws.Range("M138").Value = 2725
$ws.Range("N138").Value = -17819.3054

# ARM row 28
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 1042.75
$ws.Range("I28").Value = 1042.75
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 1042.75
$ws.Range("L28").Value = 0
$ws.Range("M28").ClearContents()
$ws.Range("N28").Value = -850.75

# ARM row 32
$ws.Range("H32").Value = 11000.52
$ws.Range("I32").Value = 9542.465
$ws.Range("J32").Value = 19957.143
$ws.Range("K32").Value = 9542.465
$ws.Range("L32").Value = 19957.143
$ws.Range("M32").Value = -9255.465
$ws.Range("N32").Value = -20531.143

# ARM row 70
$ws.Range("H70").Value = 33577
$ws.Range("J70").Value = 33577
$ws.Range("L70").Value = 33577
$ws.Range("N70").Value = -34117

# ARM row 73
$ws.Range("H73").Value = 33577
$ws.Range("J73").Value = 33577
$ws.Range("L73").Value = 33577
$ws.Range("N73").Value = -35449

# ARM row 74
$ws.Range("H74").Value = 18520770
$ws.Range("I74").Value = 33334674
$ws.Range("J74").Value = 3392
$ws.Range("K74").Value = 33334674
$ws.Range("L74").Value = 3392
$ws.Range("M74").Value = -33333800
$ws.Range("N74").Value = -5140

# ARM row 77
$ws.Range("H77").Value = 18520770
$ws.Range("I77").Value = 33334674
$ws.Range("J77").Value = 3392
$ws.Range("K77").Value = 166673370
$ws.Range("L77").Value = 16960
$ws.Range("M77").Value = -166669002
$ws.Range("N77").Value = -25696

# ARM row 99
$ws.Range("H99").Value = 1042.75
$ws.Range("I99").Value = 1042.75
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1042.75
$ws.Range("L99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = 1952.25

# ARM row 102
$ws.Range("H102").Value = 1800
$ws.Range("I102").Value = 1800
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1800
$ws.Range("L102").Value = 0
$ws.Range("M102").ClearContents()
$ws.Range("N102").Value = -178

# ARM row 117
$ws.Range("H117").Value = 37021.453
$ws.Range("J117").Value = 37021.453
$ws.Range("L117").Value = 37021.453
$ws.Range("N117").Value = -46199.453

# ARM row 122
$ws.Range("H122").Value = 5145.5757
$ws.Range("I122").Value = 5551.4287
$ws.Range("J122").Value = 2872.8
$ws.Range("K122").Value = 16654.2861
$ws.Range("L122").Value = 8618.400000000001
$ws.Range("M122").Value = -14204.2861
$ws.Range("N122").Value = -13518.4

# ARM row 132
$ws.Range("H132").Value = 4718697.5
$ws.Range("I132").Value = 6099090
$ws.Range("J132").Value = 2356.5
$ws.Range("K132").Value = 18297270
$ws.Range("L132").Value = 7069.5
$ws.Range("M132").Value = -18294740
$ws.Range("N132").Value = -12129.5

# CRP row 3
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 21333.334
$ws.Range("J3").Value = 21333.334
$ws.Range("L3").Value = 21333.334
$ws.Range("N3").Value = -21559.334

# CRP row 16
$ws.Range("H16").Value = 1920.6666
$ws.Range("I16").Value = 1617.8
$ws.Range("J16").Value = 2299.25
$ws.Range("K16").Value = 1617.8
$ws.Range("L16").Value = 2299.25
$ws.Range("M16").Value = -1330.8
$ws.Range("N16").Value = -2873.25

# CRP row 86
$ws.Range("H86").Value = 3056.5557
$ws.Range("I86").Value = 2875.1333
$ws.Range("J86").Value = 3283.3333
$ws.Range("K86").Value = 2875.1333
$ws.Range("L86").Value = 3283.3333
$ws.Range("M86").Value = -1752.1333
$ws.Range("N86").Value = -5529.3333

# CRP row 89
$ws.Range("H89").Value = 3056.5557
$ws.Range("I89").Value = 2875.1333
$ws.Range("J89").Value = 3283.3333
$ws.Range("K89").Value = 14375.6665
$ws.Range("L89").Value = 16416.6665
$ws.Range("M89").Value = -8759.666499999999
$ws.Range("N89").Value = -27648.6665

# CRP row 113
$ws.Range("H113").Value = 1920.6666
$ws.Range("I113").Value = 1617.8
$ws.Range("J113").Value = 2299.25
$ws.Range("K113").Value = 1617.8
$ws.Range("L113").Value = 2299.25
$ws.Range("M113").Value = 552.2
$ws.Range("N113").Value = -6639.25

# CUL row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 553.7586
$ws.Range("I5").Value = 267.56522
$ws.Range("J5").Value = 1650.8334
$ws.Range("K5").Value = 802.6956600000001
$ws.Range("L5").Value = 4952.5002
$ws.Range("M5").Value = -690.6956600000001
$ws.Range("N5").Value = -5176.5002

# CUL row 131
$ws.Range("H131").Value = 1889.8871
$ws.Range("I131").Value = 3230.6
$ws.Range("J131").Value = 1462
$ws.Range("K131").Value = 9691.799999999999
$ws.Range("L131").Value = 4386
$ws.Range("M131").Value = -4651.799999999999
$ws.Range("N131").Value = -14466

# CUL row 132
$ws.Range("H132").Value = 2899.875
$ws.Range("I132").Value = 2224.75
$ws.Range("J132").Value = 3575
$ws.Range("K132").Value = 20022.75
$ws.Range("L132").Value = 32175
$ws.Range("M132").Value = -17492.75
$ws.Range("N132").Value = -37235

# CUL row 135
$ws.Range("H135").Value = 553.7586
$ws.Range("I135").Value = 267.56522
$ws.Range("J135").Value = 1650.8334
$ws.Range("K135").Value = 2408.08698
$ws.Range("L135").Value = 14857.5006
$ws.Range("M135").Value = 126.91302
$ws.Range("N135").Value = -19927.5006

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 14710608
$ws.Range("I136").Value = 15626855
$ws.Range("J136").Value = 50652.5
$ws.Range("K136").Value = 46880565
$ws.Range("L136").Value = 151957.5
$ws.Range("M136").Value = -46878015
$ws.Range("N136").Value = -157057.5

# LTW row 140
$ws.Range("H140").Value = 53473.332
$ws.Range("J140").Value = 53473.332
$ws.Range("L140").Value = 53473.332
$ws.Range("N140").Value = -63833.332

# WVR row 64
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 15361
$ws.Range("J64").Value = 15361
$ws.Range("L64").Value = 15361
$ws.Range("N64").Value = -15857

# WVR row 67
$ws.Range("H67").Value = 15361
$ws.Range("J67").Value = 15361
$ws.Range("L67").Value = 15361
$ws.Range("N67").Value = -17077

# WVR row 132
$ws.Range("H132").Value = 1667.2667
$ws.Range("I132").Value = 1767.8422
$ws.Range("J132").Value = 1493.5454
$ws.Range("K132").Value = 5303.5266
$ws.Range("L132").Value = 4480.6362
$ws.Range("M132").Value = -2773.5266
$ws.Range("N132").Value = -9540.636200000001

# WVR row 135
$ws.Range("H135").Value = 29541.25
$ws.Range("J135").Value = 29541.25
$ws.Range("L135").Value = 29541.25
$ws.Range("N135").Value = -39681.25
